# "Updated items, fixed a few last minute bugs"
#
# The Affixes sheet lists prefix/suffix item modifiers. This update lowers the
# stat-modifier percentages for every tier (both the base "prefix" row and its
# five single-stat "suffix" variants), and removes the now-unused can_drop (Q)
# flag from the tiers that used to have it set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tier 1 - Balanced Energies (row 2) / single-stat suffixes (rows 3-7): 0.1/0.15 -> 0.03/0.05
$ws.Range("C2:I2").Value = 0.03
$ws.Range("F3").Value = 0.05
$ws.Range("H4").Value = 0.05
$ws.Range("G5").Value = 0.05
$ws.Range("J6").Value = 0.05
$ws.Range("I7").Value = 0.05

# Tier 2 - Natures Balancing Bliss (row 8) / suffixes (rows 9-13): 0.22/0.3 -> 0.05/0.07
$ws.Range("C8:J8").Value = 0.05
$ws.Range("Q8").ClearContents()
$ws.Range("F9").Value = 0.07
$ws.Range("Q9").ClearContents()
$ws.Range("H10").Value = 0.07
$ws.Range("Q10").ClearContents()
$ws.Range("G11").Value = 0.07
$ws.Range("Q11").ClearContents()
$ws.Range("J12").Value = 0.07
$ws.Range("Q12").ClearContents()
$ws.Range("I13").Value = 0.07
$ws.Range("Q13").ClearContents()

# Tier 3 - Queens Blessing (row 14) / suffixes (rows 15-19): 0.24/0.35 -> 0.08/0.1
$ws.Range("C14:J14").Value = 0.08
$ws.Range("Q14").ClearContents()
$ws.Range("F15").Value = 0.1
$ws.Range("Q15").ClearContents()
$ws.Range("H16").Value = 0.1
$ws.Range("Q16").ClearContents()
$ws.Range("G17").Value = 0.1
$ws.Range("Q17").ClearContents()
$ws.Range("J18").Value = 0.1
$ws.Range("Q18").ClearContents()
$ws.Range("I19").Value = 0.1
$ws.Range("Q19").ClearContents()

# Tier 4 - Wishing Spell (row 20) / suffixes (rows 21-25): 0.28/0.38 -> 0.1/0.11
$ws.Range("C20:J20").Value = 0.1
$ws.Range("F21").Value = 0.11
$ws.Range("H22").Value = 0.11
$ws.Range("G23").Value = 0.11
$ws.Range("J24").Value = 0.11
$ws.Range("I25").Value = 0.11

# Tier 5 - Chakra Alignment (row 26) / suffixes (rows 27-31): 0.34/0.44 -> 0.12/0.15
$ws.Range("C26:J26").Value = 0.12
$ws.Range("F27").Value = 0.15
$ws.Range("H28").Value = 0.15
$ws.Range("G29").Value = 0.15
$ws.Range("J30").Value = 0.15
$ws.Range("I31").Value = 0.15

# Tier 6 - Earth Tuned (row 32) / suffixes (rows 33-37): 0.38/0.48 -> 0.15/0.17
$ws.Range("C32:J32").Value = 0.15
$ws.Range("F33").Value = 0.17
$ws.Range("H34").Value = 0.17
$ws.Range("G35").Value = 0.17
$ws.Range("J36").Value = 0.17
$ws.Range("I37").Value = 0.17
